$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new changelog entry in the row right after the last used row,
# copying the formatting from the row above so the new row matches the
# existing changelog entries
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Fixed Rabbit CEM"

# Move the active selection to the next empty row, mirroring Excel's
# typical post-edit cursor placement
$ws.Range("A11").Select()
